$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "We don’t expect any risks to you if you are interviewed. If any questions make you uncomfortable, you don’t have to answer them. If you become upset during the interview, you can let your interviewer know. Remember, you can stop participating anytime without giving a reason. We care about your well-being." "Asilindelanga nayiphi na imingcipheko kuwe ukuba unodliwano-ndlebe. Ukuba nayiphi na imibuzo ikwenza ungakhululeki, akunyanzelekanga uyiphendule. Ukuba uyacaphuka ngexesha lodliwano-ndlebe, unokwazisa umenzi wodliwano-ndlebe lwakho. Khumbula, ungayeka ukuthatha inxaxheba nanini na ngaphandle kokunikeza isizathu. Siyakhathala ngempilo-ntle yakho."

Replace-Text "We also want to make sure you're safe. If we notice that you or your family are in serious danger, we might refer you for support or could need to ask for help from other places outside of this study, like social or medical services." "Sifuna nokuqinisekisa ukuba ukhuselekile. Ukuba siyaqaphelisisa ukuba wena okanye usapho lwakho lukweyona inkulu ingozi, singanithumela ukuze nifumane inkxaso okanye singadinga ukucela uncedo kwezinye indawo ezingaphandle koluphononongo, njengoo nontlalontle okanye uncedo lwezempilo."

Replace-Text "This study is part of the Global Parenting Initiative, funded by the LEGO Foundation, Oak Foundation, the World Childhood Foundation, The Human Safety Net, and the UK Research and Innovation Global Challenges Research Fund. " "Olu phononongo luyinxalenye ye Global Parenting Initiative, luxhaswe ngokwezimali ngu LEGO Foundation, Oak Fundation, i-World Childhood Foundation, i-Human Safety Net kunye ne UK Research kunye ne Innovaion Global Challenges Research Fund. "

Replace-Text "Data protection" "Ukhuseleko lwedatha"

Replace-Text "The University Cape Town makes sure your personal information is used safely and correctly, just for research. The study follows data protection laws like GDPR (General Data Protection Regulation) in the UK and POPIA (Protection of Personal Information Act) in South Africa. Any data that is transferred across borders will comply with POPIA. " "Idyunivesithi yaseKapa iqinisekisa ukuba iinkcukacha zakho zobuqu zisetyenziswa ngokukhuselekileyo nangokuchanekileyo, nje kuphando kuphela. Uphononongo lulandela imithetho yokukhuselwa kwedatha efana ne-GDPR (General Data Protection Regulation) e-UK kunye ne-POPIA (uMthetho woKhuselo loLwazi loMntu) eMzantsi Afrika. Nayiphi na idatha ethi ithunyelwe ngaphesheya kwemida izakuthobelana ne POPIA. "

Replace-Text "Who has approved this study?" "Ngubani ogunyazise oluphononongo?"

Replace-Text "[Once the ethics has been approved this will read as follows: This study has received approval from the University of Cape Town’s Centre for Social Science Research Ethics Committee and University of Cape Town’s Faculty of Health Sciences Human Research Ethics Committee. The study has also been approved by the Western Cape Department of Health and Wellness and Department of Social Development, and City of Cape Town’s City health.]" "[Yakube ivunyiwe imigaqo yokuziphatha iya kufundeka ngoluhlobo lulandelayo: Olu phononongo lufumene imvume kwiDyunivesithi yaseKapa kwiZiko leKomiti yeeNqoba zoPhando kwiNzululwazi yezeNtlalo kunye neKomiti yeeNqoba zokuziphatha zoPhando lweDyunivesithi yaseKapa. Olu phononongo lukwavunyiwe liSebe lezeMpilo leNtshona Koloni kunye ne-Mpilo kunye neSebe loPhuhliso loLuntu, kunye nesebe lempilo yeSixeko saseKapa.]"

Replace-Text "Who do I contact if I have questions or concerns?" "Ngubani endinokuqhagamshelana naye ukuba ndinemibuzo okanye iinkxalabo?"

Replace-Text "If you have any questions or concerns about your rights as a study participant, you can contact the study team at swift@globalparenting.org or on WhatsApp at +27 XX XXX XXXX (messages only)." "Ukuba unayo nayiphi na imibuzo okanye iinkxalabo malunga namalungelo akho njengomthathi-nxaxheba kuphando, ungaqhagamshelana neqela lophononongo ku-swift@globalparenting.org okanye ku-WhatsApp ku- +27 XX XXX XXXX (imiyalezo kuphela)."

Replace-Text "If you have more questions or concerns about your rights, you can contact one of the ethics committees listed: " "Ukuba uneminye imibuzo okanye iinkxalabo malunga namalungelo akho, ungaqhagamshelana nenye yee komiti yokuziphatha edwelisiweyo: "

Replace-Text "Name" "Igama"
Replace-Text "Telephone" "Inombolo yomnxeba"
Replace-Text "Email" "Imeyile"

Replace-Text "Informed Telephonic consent to take part in the study." "Imvume yoMnxeba echaziweyo yokuthatha inxaxheba kuphononongo."

# This replacement contains literal straight double quotes ("ndiyavuma"); using
# Find.Execute's replacement text here would trigger Word's smart-quote
# autocorrect and turn them into curly quotes. Locate the matching range via
# Find (search-only, no replace) on a duplicated range, then set its .Text
# directly, which performs a literal (non-autocorrected) text assignment.
$rng = $d.Content.Duplicate
$found = $rng.Find.Execute("Please respond with the word “agree” to each as I go through each of the following points. If you don’t agree, we can go over any other information you need to make your decision and if you still agree then we can proceed:", `
                            $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Text = 'Nceda uphendule ngegama elithi "ndiyavuma" kwinqaku ngalinye njengoko ndihamba ngenqaku ngalinye kulana alandelayo. Ukuba awuvumi, singajongisisa kulo naluphi na ulwazi oludingayo ukuze uthathe isigqibo kwaye ukuba usavuma singaqhubeka:'
}

Replace-Text "Someone from the research team has gone over all the information above and I know what I need to do." "Umntu osuka kwiqela lophando uye wajongisisa lonke ulwazi olungentla kwaye ndiyayazi into ekufuneka ndiyenzile."

Replace-Text "I had time to think about the information and ask questions. I am happy with the answers which I got. " "Ndiye ndanexesha lokucinga malunga nolwazi kunye nokubuza imibuzo. Ndiyavuya ngeempendulo endizifumeneyo. "

Replace-Text "I know I can say yes or no to being in the study. If I say yes, I can stop any time before the [*date to still be determined] without saying why, and nothing bad will happen." "Ndiyazi ndingathi ewe okanye hayi ekubeni yinxalenye yophononongo. Ukuba ndithi ewe, ndingayimisa nangaliphi na ixesha phambi kwe [*umhla usamiselwa] ndingatsho ukuba kutheni, kwaye akukho nto imbi iya kwenzeka."

Replace-Text "I know who can see my information after the interview, how it will be kept safe, and what happens to it after the study." "Ndiyazi ukuba ngubani onokubona ulwazi lwam emva kodliwano-ndlebe, ukuba luya kugcinwa njani lukhuselekile, kwaye kwenzeka ntoni kulo emva kophononongo."

Replace-Text "I know I can request access to my data, correct any mistakes, ask to delete it, or for it to be transferred somewhere else." "Ndiyazi ukuba ndingacela ukufikelela kwidatha yam, ndilungise naziphi na iimpazamo, ndicele ukuyicima, okanye ukuba idluliselwe kwenye indawo."

Replace-Text "I know that I won’t be named in any papers or reports from this study." "Ndiyazi ukuba andizukuchazwa kuwo nawaphi na amaphepha okanye iingxelo zolu phononongo."

Replace-Text "I know who to tell if I have a problem with the study." "Ndiyazi ukuba mandixelele bani ukuba ndinengxaki ngoluphononongo."

Replace-Text "I can be contacted again if more information is needed from me." "Ndingaqhagamshelwa kwakhona ukuba ulwazi oluninzi luyafuneka kum."

Replace-Text "I understand the team will keep my contact information safe so they can tell me about the results of the study." "Ndiyayiqonda ukuba iqela liya kugcina iinkcukacha zam zoqhagamshelwano zikhuselekile ukuze bandixelele ngeziphumo zophononongo."
